$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-seed the "RMS" shared string before "Variance" so the shared-string
# table ends up in the same order as the target workbook.
$ws.Range("A15").Value = "RMS"

# Insert a new column C ("Variance" = STD^2); this shifts the old
# Integral/Time columns from C/D to D/E and updates existing formulas'
# references automatically.
$ws.Columns.Item(3).Insert()

$ws.Range("C1").Value = "Variance"
$ws.Range("C2:C11").Formula = "=B2^2"

# Row 13 ("Avg"): average of the new Variance column.
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"

# Row 14 ("STD"): no longer computed for the STD/Variance columns, only
# for Integral (now D) and Time (now E) - those formulas already shifted
# correctly from the column insert, so just clear the old B14 STD value.
$ws.Range("B14").ClearContents()

# New row 15: RMS = sqrt(avg variance).
$ws.Range("B15").Formula = "=SQRT(C13)"

# Cosmetic: match the new column width + the workbook's final selection.
$ws.Columns.Item(3).ColumnWidth = 11.33
[void]$ws.Range("B16").Select()
